$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.331.50'
$ws.Range("E2").Value = '  -1.71%  '
$ws.Range("D3").Value = '2.454.05'
$ws.Range("E3").Value = '  -1.71%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.42'
$ws.Range("E5").Value = '  -2.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.38'
$ws.Range("E6").Value = '  -1.82%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -1.76%  '
$ws.Range("D9").Value = '2.453.60'
$ws.Range("E9").Value = '  -1.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.151'
$ws.Range("E10").Value = '  -5.34%  '
$ws.Range("E11").Value = '  -1.85%  '
$ws.Range("E12").Value = '  -5.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.80'
$ws.Range("E13").Value = '  -2.55%  '
$ws.Range("D14").Value = '2.909.38'
$ws.Range("E14").Value = '  -1.53%  '
$ws.Range("D15").Value = '68.228.69'
$ws.Range("E15").Value = '  -1.67%  '
$ws.Range("E16").Value = '  -3.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '23.29'
$ws.Range("E17").Value = '  -5.65%  '
$ws.Range("D18").Value = '2.485.42'
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.98'
$ws.Range("E19").Value = '  -1.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.19'
$ws.Range("E20").Value = '  -3.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '342.49'
$ws.Range("E21").Value = '  -1.42%  '
$ws.Range("E22").Value = '  -3.19%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").Value = '  -4.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.59'
$ws.Range("E25").Value = '  -4.39%  '
$ws.Range("E26").Value = '  +7.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.71'
$ws.Range("E27").Value = '  -5.89%  '
$ws.Range("D28").Value = '2.571.17'
$ws.Range("E28").Value = '  -2.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.15'
$ws.Range("E29").Value = '  -6.59%  '
$ws.Range("D30").Value = '0.0₃0834'
$ws.Range("E30").Value = '  -6.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.26'
$ws.Range("E31").Value = '  -7.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.40'
$ws.Range("E32").Value = '  +131.39%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '433.52'
$ws.Range("E33").Value = '  -5.02%  '
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("E35").Value = '  -3.30%  '
$ws.Range("E36").Value = '  -3.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.22'
$ws.Range("E37").Value = '  +0.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.99'
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.110'
$ws.Range("E40").Value = '  -5.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.85'
$ws.Range("E41").Value = '  -2.98%  '
$ws.Range("E42").Value = '  -3.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.45'
$ws.Range("E43").Value = '  -4.69%  '
$ws.Range("E44").Value = '  -5.02%  '
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("E46").Value = '  -6.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '133.73'
$ws.Range("E47").Value = '  -5.33%  '
$ws.Range("E48").Value = '  -3.46%  '
$ws.Range("E49").Value = '  -1.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.484'
$ws.Range("E50").Value = '  -6.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.561'
$ws.Range("E51").Value = '  -2.72%  '
